$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "Cost" column (F), shifting Cost etc. one column right
$ws.Columns("F:F").Insert()

# Set the new column header
$ws.Range("F1").Value2 = "ISBN13"

# Fill in ISBN13 values for each book row
$isbns = @(9781780226583, 9780170364379, 9780684838281, 9780415583367, 9780316176200, 9780415690157, 9780140254037, 9780415623230, 9781571107282, 9781137278814)
for ($i = 0; $i -lt $isbns.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 6).Value2 = $isbns[$i]
}

# Apply number formatting (integer, no thousands separator) to the ISBN13 column
$ws.Range("F1:F11").NumberFormat = "0"

# Header cell keeps bold font (Calibri) like other header cells, F2:F11 use Arial 11 font colored FF333333
$ws.Range("F2:F11").Font.Name = "Arial"
$ws.Range("F2:F11").Font.Size = 11
$ws.Range("F2:F11").Font.Color = 3355443

# Column F formatting: best-fit width
$ws.Columns("F:F").AutoFit()

# Update the selected cell, mirroring the diff's final cursor position
$ws.Range("D14").Select()
